$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 1152-1153, pushing the existing rows 1152-1174 down to 1154-1176.
$ws.Rows("1152:1153").Insert()

# --- New row 1152 : Alcachofa, Española, Primera, $/unidad ---
$ws.Cells.Item(1152, 1).Value = 6
$ws.Cells.Item(1152, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1152, 3).Value = "Metropolitana"
$ws.Cells.Item(1152, 4).Value = 45239
$ws.Cells.Item(1152, 5).Value = 13
$ws.Cells.Item(1152, 6).Value = 100112013
$ws.Cells.Item(1152, 7).Value = "Alcachofa"
$ws.Cells.Item(1152, 8).Value = "Española"
$ws.Cells.Item(1152, 9).Value = "Primera"
$ws.Cells.Item(1152, 10).Value = 6700
$ws.Cells.Item(1152, 11).Value = 400
$ws.Cells.Item(1152, 12).Value = 450
$ws.Cells.Item(1152, 13).Value = 424
$ws.Cells.Item(1152, 14).Value = "$/unidad"
$ws.Cells.Item(1152, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1152, 16).Value = 424
$ws.Cells.Item(1152, 17).Value = 1
$ws.Cells.Item(1152, 18).Value = "Hortaliza"

# --- New row 1153 : Alcachofa, Española, Segunda, $/unidad ---
$ws.Cells.Item(1153, 1).Value = 6
$ws.Cells.Item(1153, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(1153, 3).Value = "Metropolitana"
$ws.Cells.Item(1153, 4).Value = 45239
$ws.Cells.Item(1153, 5).Value = 13
$ws.Cells.Item(1153, 6).Value = 100112013
$ws.Cells.Item(1153, 7).Value = "Alcachofa"
$ws.Cells.Item(1153, 8).Value = "Española"
$ws.Cells.Item(1153, 9).Value = "Segunda"
$ws.Cells.Item(1153, 10).Value = 5000
$ws.Cells.Item(1153, 11).Value = 300
$ws.Cells.Item(1153, 12).Value = 350
$ws.Cells.Item(1153, 13).Value = 325
$ws.Cells.Item(1153, 14).Value = "$/unidad"
$ws.Cells.Item(1153, 15).Value = "Región Metropolitana"
$ws.Cells.Item(1153, 16).Value = 325
$ws.Cells.Item(1153, 17).Value = 1
$ws.Cells.Item(1153, 18).Value = "Hortaliza"
